$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume figures with the latest scraped values.
# The cells store text (not numbers/percentages), so we force a Text number format
# before assigning, which keeps numeric-looking strings like "304.33" or "4.45%"
# from being auto-converted by Excel into numeric/percentage values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.45%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "15.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.082"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.53%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07843"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.29%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.282"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.93%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.146"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.24%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.019"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.68%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9284"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.08%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1002"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.08%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1827"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.04%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08703"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.44%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03386"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.13%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09906"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.26%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001482"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.05%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.35%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.484"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.45%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.39%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3434"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.07%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1321"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.83%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.535"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "9.26%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2237"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.39%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04663"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.36%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.40%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004493"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.49%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001297"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.18%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002695"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-20.57%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01769"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "9.72%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04703"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.96%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007828"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.20%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.13%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008435"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-14.21%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002207"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.23%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009184"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.15%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006048"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.84%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.15%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.787"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "117.99%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002685"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.30%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.15%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.15%"
